$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.689.63"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.342.97"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.02%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "3.337.40"
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "3.923.46"
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.74%  "
$ws.Range("D16").Value = "65.703.44"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000170"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "3.344.96"
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.522"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("E35").Value = "  -5.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").Value = "2.673.50"
$ws.Range("E42").Value = "  -5.55%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0668"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "336.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0280"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.971"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
Write-Host "Updated crypto prices and volumes"
